$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-06 Saturday", "2024-01-07 Sunday"),
    @("33×30=", "45×77="),
    @("59×39=", "97×99="),
    @("73×37=", "60×38="),
    @("89×34=", "15×62="),
    @("99×90=", "25×85="),
    @("79×32=", "81×58="),
    @("97×92=", "37×18="),
    @("92×40=", "53×88="),
    @("82×67=", "44×31="),
    @("98×39=", "38×43="),
    @("43×65=", "91×38="),
    @("89×62=", "27×58="),
    @("53×84=", "29×73="),
    @("17×31=", "90×50="),
    @("71×93=", "79×85="),
    @("29×87=", "46×38="),
    @("19×60=", "42×47="),
    @("99×65=", "13×73="),
    @("17×69=", "76×90="),
    @("30×18=", "93×31="),
    @("90×15=", "70×29="),
    @("70×62=", "94×11="),
    @("58×81=", "45×89="),
    @("89×49=", "65×29="),
    @("97×81=", "40×56=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
